# Update the "Date" metadata value and replace the SNOMED CT ECL-constraint
# based ValueSet.compose.include with a full code-system import from the
# EDQM Standard Terms code system ("suppression page etude complementaire").

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the Date property -------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-04-09T16:24:06+00:00"

# --- Include #0 sheet: swap the constraint filter for a full code import --
$wsInc = $wb.Worksheets.Item("Include #0")

# Row 1 used to be Property | Operation | Value -> now just "Codes"
$wsInc.Range("A1").Value = "Codes"
$wsInc.Range("B1").Clear()
$wsInc.Range("C1").Clear()

# Row 2 used to be constraint | = | <ECL expression> -> now just "All codes"
$wsInc.Range("A2").Value = "All codes"
$wsInc.Range("B2").Clear()
$wsInc.Range("C2").Clear()

# Row 4: System URI stays, but now points at the EDQM Standard Terms system
$wsInc.Range("A4").Value = "System URI"
$wsInc.Range("B4").Value = "http://standardterms.edqm.eu"
